$wb = $excel.ActiveWorkbook

# The file f0951972-e5b0-4340-9dbf-16b1f00d58a1.md has been handed off for
# localization again: bump its status to "Ready for handoff" in both the
# per-language sheets and the Overview roll-up, update its priority to "mt",
# and record the new handoff timestamps.

# --- zh-cn sheet (row 3 = f0951972 file) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("H3").Value = "2016-08-23 18:13:56"

# --- de-de sheet (row 3 = f0951972 file) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "mt"
$wsDe.Range("H3").Value = "2016-08-23 18:14:02"

# --- Overview sheet (row 3 = f0951972 file) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-23 18:14:02"
